$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.933.59"
$ws.Range("E2").Value = "  +5.14%  "
$ws.Range("D3").Value = "3.503.03"
$ws.Range("E3").Value = "  +2.74%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'593.96"
$ws.Range("E5").Value = "  +4.37%  "
$ws.Range("D6").Value = "'169.03"
$ws.Range("E6").Value = "  +7.41%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "3.505.63"
$ws.Range("E8").Value = "  +2.70%  "
$ws.Range("D9").Value = "'0.575"
$ws.Range("E9").Value = "  +1.33%  "
$ws.Range("E10").Value = "  +0.55%  "
$ws.Range("E11").Value = "  +5.40%  "
$ws.Range("D12").Value = "'0.439"
$ws.Range("E12").Value = "  +4.01%  "
$ws.Range("D13").Value = "4.108.51"
$ws.Range("E13").Value = "  +2.83%  "
$ws.Range("E14").Value = "  +0.18%  "
$ws.Range("D15").Value = "'28.20"
$ws.Range("E15").Value = "  +4.42%  "
$ws.Range("E16").Value = "  +4.33%  "
$ws.Range("D17").Value = "66.885.80"
$ws.Range("E17").Value = "  +4.99%  "
$ws.Range("D18").Value = "3.522.35"
$ws.Range("E18").Value = "  +4.26%  "
$ws.Range("D19").Value = "'6.32"
$ws.Range("E19").Value = "  +3.86%  "
$ws.Range("D20").Value = "'14.06"
$ws.Range("E20").Value = "  +3.46%  "
$ws.Range("D21").Value = "'395.00"
$ws.Range("E21").Value = "  +2.36%  "
$ws.Range("D22").Value = "'7.95"
$ws.Range("E22").Value = "  +2.43%  "
$ws.Range("D23").Value = "'73.43"
$ws.Range("E23").Value = "  +3.14%  "
$ws.Range("E24").Value = "  +11.87%  "
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("D26").Value = "'0.530"
$ws.Range("E26").Value = "  +2.97%  "
$ws.Range("D27").Value = "'10.04"
$ws.Range("E27").Value = "  +3.91%  "
$ws.Range("D28").Value = "'0.182"
$ws.Range("E28").Value = "  +2.49%  "
$ws.Range("D30").Value = "'6.37"
$ws.Range("E30").Value = "  +4.79%  "
$ws.Range("E31").Value = "  +5.89%  "
$ws.Range("D32").Value = "'2.06"
$ws.Range("E32").Value = "  +4.48%  "
$ws.Range("D33").Value = "'23.56"
$ws.Range("E33").Value = "  +3.08%  "
$ws.Range("D34").Value = "'7.44"
$ws.Range("E34").Value = "  +7.25%  "
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("D36").Value = "'1.60"
$ws.Range("E36").Value = "  +6.02%  "
$ws.Range("D37").Value = "'162.49"
$ws.Range("E37").Value = "  +1.16%  "
$ws.Range("D38").Value = "'0.900"
$ws.Range("E38").Value = "  +6.89%  "
$ws.Range("E39").Value = "  +6.58%  "
$ws.Range("E40").Value = "  +4.12%  "
$ws.Range("D41").Value = "'4.66"
$ws.Range("E41").Value = "  +7.27%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "'6.68"
$ws.Range("E42").Value = "  +5.14%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.836.55"
$ws.Range("E43").Value = "  +1.89%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'26.36"
$ws.Range("E44").Value = "  +1.84%  "
$ws.Range("D45").Value = "'43.50"
$ws.Range("E45").Value = "  +1.20%  "
$ws.Range("D46").Value = "'26.58"
$ws.Range("E46").Value = "  +2.63%  "
$ws.Range("D47").Value = "'0.0316"
$ws.Range("E47").Value = "  +4.39%  "
$ws.Range("D48").Value = "'2.55"
$ws.Range("E48").Value = "  +7.63%  "
$ws.Range("D49").Value = "'348.44"
$ws.Range("E49").Value = "  +6.53%  "
$ws.Range("D50").Value = "'1.08"
$ws.Range("E50").Value = "  +5.09%  "
$ws.Range("D51").Value = "'33.66"
$ws.Range("E51").Value = "  +11.86%  "
